$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"0.1059284729039495"
$ws.Range("B3").Value = [double]"0.005697116712966829"
$ws.Range("C3").Value = [double]"0.0007490514149459492"
$ws.Range("D3").Value = [double]"5.121242706737477"
$ws.Range("E3").Value = [double]"0.06724021490028202"
$ws.Range("F3").Value = [double]"0.004228997636287662"
$ws.Range("G3").Value = [double]"0.007165235789645995"
$ws.Range("H3").Value = [double]"0.1116255896169164"
$ws.Range("B4").Value = [double]"0.008152729322212766"
$ws.Range("C4").Value = [double]"0.001337471587661834"
$ws.Range("D4").Value = [double]"5.495322936491093"
$ws.Range("E4").Value = [double]"0.06074208409273012"
$ws.Range("F4").Value = [double]"0.005531324789590864"
$ws.Range("G4").Value = [double]"0.01077413385483467"
$ws.Range("H4").Value = [double]"0.1140812022261623"
$ws.Range("B5").Value = [double]"0.01221126351024392"
$ws.Range("C5").Value = [double]"0.006873872428501557"
$ws.Range("D5").Value = [double]"4.762879975651003"
$ws.Range("E5").Value = [double]"0.1459366223318906"
$ws.Range("F5").Value = [double]"-0.001261331335154671"
$ws.Range("G5").Value = [double]"0.02568385835564252"
$ws.Range("H5").Value = [double]"0.1181397364141935"
$ws.Range("B6").Value = [double]"0.01225744341491095"
$ws.Range("C6").Value = [double]"0.003235812952605201"
$ws.Range("D6").Value = [double]"3.152175087274364"
$ws.Range("E6").Value = [double]"0.06848290591496634"
$ws.Range("F6").Value = [double]"0.005915348277778384"
$ws.Range("G6").Value = [double]"0.01859953855204351"
$ws.Range("H6").Value = [double]"0.1181859163188605"
$ws.Range("B7").Value = [double]"0.01136641791412235"
$ws.Range("C7").Value = [double]"0.002952467626870178"
$ws.Range("D7").Value = [double]"2.986788412109476"
$ws.Range("E7").Value = [double]"0.01082846962630187"
$ws.Range("F7").Value = [double]"0.005579670892552515"
$ws.Range("G7").Value = [double]"0.01715316493569218"
$ws.Range("H7").Value = [double]"0.1172948908180719"
$ws.Range("B8").Value = [double]"0.01208888530555888"
$ws.Range("C8").Value = [double]"0.002308071774037399"
$ws.Range("D8").Value = [double]"2.748036709434197"
$ws.Range("E8").Value = [double]"0.02662049954933978"
$ws.Range("F8").Value = [double]"0.007565134630410122"
$ws.Range("G8").Value = [double]"0.01661263598070764"
$ws.Range("H8").Value = [double]"0.1180173582095084"
$ws.Range("B9").Value = [double]"0.008387450226873684"
$ws.Range("C9").Value = [double]"0.004701658509627594"
$ws.Range("D9").Value = [double]"1.638319251714986"
$ws.Range("E9").Value = [double]"0.06007979086243981"
$ws.Range("F9").Value = [double]"-0.0008276573458938054"
$ws.Range("G9").Value = [double]"0.01760255779964117"
$ws.Range("H9").Value = [double]"0.1143159231308232"
$ws.Range("B10").Value = [double]"-0.1059284729039495"
$ws.Range("C10").Value = [double]"0.0005319173889280426"
$ws.Range("D10").Value = [double]"-221.2603094004441"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.1069710153701096"
$ws.Range("G10").Value = [double]"-0.1048859304377894"
$ws.Range("B11").Value = [double]"-0.0505401812519105"
$ws.Range("C11").Value = [double]"0.0005752973280551824"
$ws.Range("D11").Value = [double]"-95.29271873177335"
$ws.Range("E11").Value = [double]"5.301902137264934e-183"
$ws.Range("F11").Value = [double]"-0.0516677471240825"
$ws.Range("G11").Value = [double]"-0.04941261537973848"
$ws.Range("H11").Value = [double]"0.05538829165203905"
$ws.Range("B12").Value = [double]"-0.03924514760045478"
$ws.Range("C12").Value = [double]"0.0005514303451821708"
$ws.Range("D12").Value = [double]"-76.75585415038296"
$ws.Range("E12").Value = [double]"2.278889107319344e-105"
$ws.Range("F12").Value = [double]"-0.04032593489731532"
$ws.Range("G12").Value = [double]"-0.03816436030359425"
$ws.Range("H12").Value = [double]"0.06668332530349477"
$ws.Range("B13").Value = [double]"-0.03577771674960952"
$ws.Range("C13").Value = [double]"0.0005459076748464473"
$ws.Range("D13").Value = [double]"-70.20722794576045"
$ws.Range("E13").Value = [double]"5.209872167739041e-89"
$ws.Range("F13").Value = [double]"-0.03684767976462736"
$ws.Range("G13").Value = [double]"-0.03470775373459169"
$ws.Range("H13").Value = [double]"0.07015075615434002"
$ws.Range("B14").Value = [double]"-0.03137846230243427"
$ws.Range("C14").Value = [double]"0.0005367404804594845"
$ws.Range("D14").Value = [double]"-63.59836359913692"
$ws.Range("E14").Value = [double]"1.434082905331789e-16"
$ws.Range("F14").Value = [double]"-0.03243045788956587"
$ws.Range("G14").Value = [double]"-0.03032646671530268"
$ws.Range("H14").Value = [double]"0.07455001060151528"
$ws.Range("B15").Value = [double]"-0.0289721421231689"
$ws.Range("C15").Value = [double]"0.0005280113599322924"
$ws.Range("D15").Value = [double]"-58.81867679573283"
$ws.Range("E15").Value = [double]"1.579264733050728e-32"
$ws.Range("F15").Value = [double]"-0.03000702889449998"
$ws.Range("G15").Value = [double]"-0.02793725535183783"
$ws.Range("H15").Value = [double]"0.07695633078078065"
$ws.Range("B16").Value = [double]"-0.02750203918285931"
$ws.Range("C16").Value = [double]"0.0005251949708560035"
$ws.Range("D16").Value = [double]"-56.25185481906498"
$ws.Range("E16").Value = [double]"9.186402907995585e-14"
$ws.Range("F16").Value = [double]"-0.02853140591240675"
$ws.Range("G16").Value = [double]"-0.02647267245331187"
$ws.Range("H16").Value = [double]"0.07842643372109025"
$ws.Range("B17").Value = [double]"-0.02546824286065745"
$ws.Range("C17").Value = [double]"0.0005342416321700299"
$ws.Range("D17").Value = [double]"-51.1574066992342"
$ws.Range("E17").Value = [double]"2.936017371591967e-07"
$ws.Range("F17").Value = [double]"-0.02651534078211408"
$ws.Range("G17").Value = [double]"-0.0244211449392008"
$ws.Range("H17").Value = [double]"0.0804602300432921"
$ws.Range("B18").Value = [double]"-0.02262277342282905"
$ws.Range("C18").Value = [double]"0.0005329701905898541"
$ws.Range("D18").Value = [double]"-44.94770975042273"
$ws.Range("E18").Value = [double]"0.03244842495338483"
$ws.Range("F18").Value = [double]"-0.02366737935363432"
$ws.Range("G18").Value = [double]"-0.02157816749202379"
$ws.Range("H18").Value = [double]"0.0833056994811205"
$ws.Range("B19").Value = [double]"-0.01892621495415249"
$ws.Range("C19").Value = [double]"0.0005276324108027382"
$ws.Range("D19").Value = [double]"-37.24444043834239"
$ws.Range("E19").Value = [double]"0.0003565824974255133"
$ws.Range("F19").Value = [double]"-0.01996035900883707"
$ws.Range("G19").Value = [double]"-0.01789207089946792"
$ws.Range("H19").Value = [double]"0.08700225794979706"
$ws.Range("B20").Value = [double]"-0.01570785091646774"
$ws.Range("C20").Value = [double]"0.0005430917219379897"
$ws.Range("D20").Value = [double]"-28.81888216264499"
$ws.Range("E20").Value = [double]"0.05874676758181339"
$ws.Range("F20").Value = [double]"-0.01677229477380077"
$ws.Range("G20").Value = [double]"-0.0146434070591347"
$ws.Range("H20").Value = [double]"0.09022062198748182"
$ws.Range("B21").Value = [double]"-0.0136455070299031"
$ws.Range("C21").Value = [double]"0.0005521607385374306"
$ws.Range("D21").Value = [double]"-22.90470114236884"
$ws.Range("E21").Value = [double]"0.09340033721905226"
$ws.Range("F21").Value = [double]"-0.01472772588841303"
$ws.Range("G21").Value = [double]"-0.01256328817139317"
$ws.Range("H21").Value = [double]"0.09228296587404645"
$ws.Range("B22").Value = [double]"-0.01144997290304208"
$ws.Range("C22").Value = [double]"0.0005508601744417244"
$ws.Range("D22").Value = [double]"-18.51328666375814"
$ws.Range("E22").Value = [double]"0.06127779573725212"
$ws.Range("F22").Value = [double]"-0.01252964268663251"
$ws.Range("G22").Value = [double]"-0.01037030311945164"
$ws.Range("H22").Value = [double]"0.09447850000090748"
$ws.Range("B23").Value = [double]"-0.01021020149541067"
$ws.Range("C23").Value = [double]"0.0005569828254148677"
$ws.Range("D23").Value = [double]"-16.22444714055327"
$ws.Range("E23").Value = [double]"0.05607413366788849"
$ws.Range("F23").Value = [double]"-0.01130187152770395"
$ws.Range("G23").Value = [double]"-0.009118531463117391"
$ws.Range("H23").Value = [double]"0.09571827140853888"
$ws.Range("B24").Value = [double]"-0.008718534987623039"
$ws.Range("C24").Value = [double]"0.000552342838046243"
$ws.Range("D24").Value = [double]"-12.9600644533773"
$ws.Range("E24").Value = [double]"0.08841046935965835"
$ws.Range("F24").Value = [double]"-0.009801110764317524"
$ws.Range("G24").Value = [double]"-0.007635959210928555"
$ws.Range("H24").Value = [double]"0.09720993791632651"
$ws.Range("B25").Value = [double]"-0.005390714992994997"
$ws.Range("C25").Value = [double]"0.0005480546881707046"
$ws.Range("D25").Value = [double]"-6.897097142683885"
$ws.Range("E25").Value = [double]"0.02330089830328605"
$ws.Range("F25").Value = [double]"-0.006464886107032022"
$ws.Range("G25").Value = [double]"-0.004316543878957971"
$ws.Range("H25").Value = [double]"0.1005377579109546"
$ws.Range("B26").Value = [double]"0.01616126310265665"
$ws.Range("C26").Value = [double]"0.005246430916884039"
$ws.Range("D26").Value = [double]"9.377127405638342"
$ws.Range("E26").Value = [double]"0.1132196034960564"
$ws.Range("F26").Value = [double]"0.005878419296894312"
$ws.Range("G26").Value = [double]"0.02644410690841898"
$ws.Range("H26").Value = [double]"0.1220897360066062"
